# "updated severity for parameters"
# Sheet1: rows 8 and 9 (E8, E9) had their Severity changed from "High" to "Low".
# Also the sheet's scroll position / active selection moved (topLeftCell A4,
# active cell E10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E8").Value = "Low"
$ws.Range("E9").Value = "Low"

# Make sure Sheet1 is the active sheet/window and restore the view state:
# scrolled so row 4 is at the top, with E10 as the active/selected cell.
$ws.Activate()
$ws.Range("E10").Select()
$excel.ActiveWindow.ScrollRow = 4
